# Adicionando respostas do quiz da semana 5
#
# Colors every answer paragraph of the quiz (the multiple-choice options
# and the free-response answers) green (RGB 00FF00) to mark them as the
# correct / filled-in answers, leaving the question text untouched.
#
# wdColor value 65280 (0x00FF00) == RGB(0, 255, 0) -> OOXML <w:color w:val="00FF00"/>

$d = $word.ActiveDocument
$green = 65280

# 1-based paragraph indexes of every answer run that needs to turn green.
$answerParagraphs = @(2, 3, 4, 5, 8, 9, 10, 11, 14, 17, 20, 21, 22, 23, 24, 27, 28, 29, 30, 31, 34, 35, 38, 39, 42, 43)

foreach ($idx in $answerParagraphs) {
    $p = $d.Paragraphs($idx)
    $pRange = $p.Range
    # Re-wrap the paragraph's start/end as a fresh Range so the color change
    # lands only on the run(s) inside the paragraph, not on the paragraph
    # mark's own run properties (w:pPr/w:rPr).
    $rng = $d.Range($pRange.Start, $pRange.End)
    $rng.Font.Color = $green
}
